$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update revised values for existing rows 785-799 (columns F and G) ---
$ws.Range("F785").Value = 7179

$ws.Range("F786").Value = 6337

$ws.Range("F787").Value = 2275
$ws.Range("G787").Value = 96

$ws.Range("F788").Value = 1762

$ws.Range("F789").Value = 7874

$ws.Range("F790").Value = 4656

$ws.Range("F791").Value = 4237

$ws.Range("F792").Value = 3910

$ws.Range("F793").Value = 3394

$ws.Range("F794").Value = 1291

$ws.Range("F795").Value = 1158

$ws.Range("F796").Value = 4629
$ws.Range("G796").Value = 260

$ws.Range("F797").Value = 3389
$ws.Range("G797").Value = 166

$ws.Range("F798").Value = 3495
$ws.Range("G798").Value = 136

$ws.Range("F799").Value = 3059
$ws.Range("G799").Value = 125

# --- Append new daily rows 800-809 ---
$newRows = @(
    @(44694, 1786914, 2132, 335, 20028, 3050, 110),
    @(44695, 1787093, 1207, 179, 20041, 1231, 57),
    @(44696, 1787142, 439, 49, 20050, 962, 45),
    @(44697, 1787404, 2275, 262, 20062, 3737, 170),
    @(44698, 1787657, 1933, 253, 20069, 2593, 79),
    @(44699, 1787919, 1746, 262, 20073, 2272, 80),
    @(44700, 1788133, 1524, 214, 20075, 2677, 49),
    @(44701, 1788334, 1684, 201, 20077, 1896, 60),
    @(44702, 1788452, 890, 118, 20080, 512, 19),
    @(44703, 1788490, 404, 38, 20081, 312, 12)
)

$rowIndex = 800
foreach ($rowData in $newRows) {
    $ws.Cells.Item($rowIndex, 1).Value = $rowData[0]
    $ws.Cells.Item($rowIndex, 2).Value = $rowData[1]
    $ws.Cells.Item($rowIndex, 3).Value = $rowData[2]
    $ws.Cells.Item($rowIndex, 4).Value = $rowData[3]
    $ws.Cells.Item($rowIndex, 5).Value = $rowData[4]
    $ws.Cells.Item($rowIndex, 6).Value = $rowData[5]
    $ws.Cells.Item($rowIndex, 7).Value = $rowData[6]
    $rowIndex++
}
